$wb = $excel.ActiveWorkbook

# --- Sheet "Metadata" (sheet1) ---
$ws = $wb.Worksheets.Item("Metadata")

# Version: 5.0.0 -> 6.0.0
$ws.Range("B3").Value = "6.0.0"

# Date: 2021-12-16T17:36:56+00:00 -> 2022-01-21T20:46:54+00:00
$ws.Range("B8").Value = "2022-01-21T20:46:54+00:00"

# Publisher value was blank -> now "Alvearie Team"
$ws.Range("B9").Value = "Alvearie Team"

# Old row 10 was "Contact" / "No display for ContactDetail" (duplicated on row 11).
# It becomes "Jurisdiction" / "United States of America", and the old duplicate
# row 11 is removed entirely (rows below shift up).
$ws.Range("A10").Value = "Jurisdiction"
$ws.Range("B10").Value = "United States of America"
$ws.Rows.Item(11).Delete()

# --- Sheet "Elements" (sheet2) ---
$ws2 = $wb.Worksheets.Item("Elements")

# Row 2 (the root "Extension" element) now carries a specific Short/Definition
# instead of the generic placeholders.
$ws2.Range("K2").Value = "Average Wholesale Price"
$ws2.Range("L2").Value = "The average wholesale price charged by wholesalers for the specific drug"
